$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking strings that must remain text
# (matches the source data which stores these as literal text, not numbers)
$ws.Range("D4,D5,D6,D7,D11,D12,D13,D16,D17,D20,D21,D22,D23,D24,D25,D26,D27,D30,D31,D32,D33,D35,D36,D37,D38,D39,D41,D43,D45,D47,D48,D50,D51").NumberFormat = "@"

$ws.Range('D2').Value = '69.114.88'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '3.500.83'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '602.92'
$ws.Range('E5').Value = '  +4.57%  '
$ws.Range('D6').Value = '170.18'
$ws.Range('E6').Value = '  -1.96%  '
$ws.Range('D7').Value = '0.610'
$ws.Range('E7').Value = '  -0.80%  '
$ws.Range('D8').Value = '3.494.91'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  +3.16%  '
$ws.Range('D11').Value = '6.76'
$ws.Range('E11').Value = '  +2.27%  '
$ws.Range('D12').Value = '0.578'
$ws.Range('E12').Value = '  -3.58%  '
$ws.Range('D13').Value = '47.17'
$ws.Range('E13').Value = '  +0.17%  '
$ws.Range('E14').Value = '  +1.59%  '
$ws.Range('D15').Value = '4.065.50'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('B16').Value = 'BitcoinCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D16').Value = '617.12'
$ws.Range('E16').Value = '  -9.10%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = '8.36'
$ws.Range('E17').Value = '  -5.66%  '
$ws.Range('D18').Value = '3.511.47'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').Value = '69.203.48'
$ws.Range('E19').Value = '  +0.54%  '
$ws.Range('D20').Value = '0.119'
$ws.Range('E20').Value = '  -2.07%  '
$ws.Range('D21').Value = '17.22'
$ws.Range('E21').Value = '  -1.41%  '
$ws.Range('D22').Value = '11.19'
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D23').Value = '0.876'
$ws.Range('E23').Value = '  -2.84%  '
$ws.Range('D24').Value = '15.80'
$ws.Range('E24').Value = '  -2.95%  '
$ws.Range('D25').Value = '96.11'
$ws.Range('E25').Value = '  -1.22%  '
$ws.Range('D26').Value = '3.85'
$ws.Range('E26').Value = '  +0.40%  '
$ws.Range('D27').Value = '5.89'
$ws.Range('E27').Value = '  +2.87%  '
$ws.Range('E28').Value = '  +0.20%  '
$ws.Range('E29').Value = '  -1.17%  '
$ws.Range('D30').Value = '9.19'
$ws.Range('E30').Value = '  -1.70%  '
$ws.Range('D31').Value = '33.32'
$ws.Range('E31').Value = '  +1.08%  '
$ws.Range('D32').Value = '8.44'
$ws.Range('E32').Value = '  -3.73%  '
$ws.Range('D33').Value = '3.10'
$ws.Range('E33').Value = '  -1.66%  '
$ws.Range('E34').Value = '  -2.12%  '
$ws.Range('D35').Value = '6.88'
$ws.Range('E35').Value = '  -4.81%  '
$ws.Range('D36').Value = '570.61'
$ws.Range('E36').Value = '  +1.61%  '
$ws.Range('D37').Value = '10.74'
$ws.Range('E37').Value = '  -0.90%  '
$ws.Range('D38').Value = '3.53'
$ws.Range('E38').Value = '  -2.09%  '
$ws.Range('D39').Value = '57.10'
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('E40').Value = '  -3.59%  '
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('E42').Value = '  +0.37%  '
$ws.Range('D43').Value = '0.0440'
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').Value = '3.394.25'
$ws.Range('E44').Value = '  -1.56%  '
$ws.Range('D45').Value = '0.325'
$ws.Range('E45').Value = '  -2.81%  '
$ws.Range('D46').Value = '0.0₃0710'
$ws.Range('E46').Value = '  +1.50%  '
$ws.Range('D47').Value = '32.75'
$ws.Range('E47').Value = '  -1.74%  '
$ws.Range('D48').Value = '2.57'
$ws.Range('E48').Value = '  -0.58%  '
$ws.Range('E49').Value = '  -2.27%  '
$ws.Range('D50').Value = '0.129'
$ws.Range('E50').Value = '  -3.36%  '
$ws.Range('D51').Value = '134.16'
$ws.Range('E51').Value = '  -0.25%  '
